$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 217.33333
$ws.Range("I9").Value = 176.33333
$ws.Range("K9").Value = 176.33333
$ws.Range("M9").Value = -7.333329999999989

$ws.Range("H12").Value = 156
$ws.Range("I12").Value = 165.42857
$ws.Range("J12").Value = 90
$ws.Range("K12").Value = 165.42857
$ws.Range("L12").Value = 90
$ws.Range("M12").Value = 4.571429999999992
$ws.Range("N12").Value = -430

$ws.Range("H15").Value = 1748.84
$ws.Range("I15").Value = 1748.84
$ws.Range("K15").Value = 5246.52
$ws.Range("M15").Value = -5077.52

$ws.Range("H33").Value = 2246.4119
$ws.Range("I33").Value = 2177.2593
$ws.Range("K33").Value = 2177.2593
$ws.Range("M33").Value = -1948.2593

$ws.Range("H69").Value = 40003.75
$ws.Range("J69").Value = 40003.75
$ws.Range("L69").Value = 120011.25
$ws.Range("N69").Value = -121759.25

$ws.Range("H72").Value = 40003.75
$ws.Range("J72").Value = 40003.75
$ws.Range("L72").Value = 360033.75
$ws.Range("N72").Value = -368769.75

$ws.Range("H107").Value = 4072.8572
$ws.Range("J107").Value = 4084.1667
$ws.Range("L107").Value = 4084.1667
$ws.Range("N107").Value = -7924.1667

$ws.Range("H132").Value = 2397.4736
$ws.Range("I132").Value = 2078.9395
$ws.Range("J132").Value = 4499.8
$ws.Range("K132").Value = 6236.818499999999
$ws.Range("L132").Value = 13499.4
$ws.Range("M132").Value = -3706.818499999999
$ws.Range("N132").Value = -18559.4

$ws.Range("H135").Value = 9315.200000000001
$ws.Range("I135").Value = 1792.9
$ws.Range("K135").Value = 16136.1
$ws.Range("M135").Value = -13601.1

$ws.Range("H138").Value = 2688.5833
$ws.Range("I138").Value = 683.6316
$ws.Range("J138").Value = 3274.6462
$ws.Range("K138").Value = 2050.8948
$ws.Range("L138").Value = 9823.938600000001
$ws.Range("M138").Value = 3089.1052
$ws.Range("N138").Value = -20103.9386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2011.8889
$ws.Range("I2").Value = 2062
$ws.Range("J2").Value = 1611
$ws.Range("K2").Value = 2062
$ws.Range("L2").Value = 1611
$ws.Range("M2").Value = -1949
$ws.Range("N2").Value = -1837

$ws.Range("H45").Value = 2381.818
$ws.Range("I45").Value = 1200
$ws.Range("K45").Value = 1200
$ws.Range("M45").Value = -823

$ws.Range("H116").Value = 2011.8889
$ws.Range("I116").Value = 2062
$ws.Range("J116").Value = 1611
$ws.Range("K116").Value = 2062
$ws.Range("L116").Value = 1611
$ws.Range("M116").Value = 232
$ws.Range("N116").Value = -6199

$ws.Range("H122").Value = 2236.25
$ws.Range("I122").Value = 1984.2858
$ws.Range("K122").Value = 5952.857400000001
$ws.Range("M122").Value = -3502.857400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2011.8889
$ws.Range("I3").Value = 2062
$ws.Range("J3").Value = 1611
$ws.Range("K3").Value = 2062
$ws.Range("L3").Value = 1611
$ws.Range("M3").Value = -1948
$ws.Range("N3").Value = -1839

$ws.Range("H107").Value = 2291.1
$ws.Range("I107").Value = 1986.1428
$ws.Range("K107").Value = 1986.1428
$ws.Range("M107").Value = -66.14280000000008

$ws.Range("H122").Value = 74995
$ws.Range("J122").Value = 74995
$ws.Range("L122").Value = 74995
$ws.Range("N122").Value = -84795

$ws.Range("H125").Value = 104495
$ws.Range("J125").Value = 104495
$ws.Range("L125").Value = 104495
$ws.Range("N125").Value = -114335

$ws.Range("H126").Value = 88995
$ws.Range("J126").Value = 88995
$ws.Range("L126").Value = 88995
$ws.Range("N126").Value = -98875

$ws.Range("H134").Value = 40214.152
$ws.Range("I134").Value = 1720.4783
$ws.Range("J134").Value = 335332.34
$ws.Range("K134").Value = 5161.4349
$ws.Range("L134").Value = 1005997.02
$ws.Range("M134").Value = -2626.4349
$ws.Range("N134").Value = -1011067.02

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 5700.8335
$ws.Range("I7").Value = 345.25
$ws.Range("J7").Value = 16412
$ws.Range("K7").Value = 345.25
$ws.Range("L7").Value = 16412
$ws.Range("M7").Value = -232.25
$ws.Range("N7").Value = -16638

$ws.Range("H22").Value = 499.5
$ws.Range("I22").Value = 499.5
$ws.Range("K22").Value = 499.5
$ws.Range("M22").Value = -149.5

$ws.Range("H88").Value = 32460.5
$ws.Range("J88").Value = 32460.5
$ws.Range("L88").Value = 32460.5
$ws.Range("N88").Value = -33272.5

$ws.Range("H91").Value = 32460.5
$ws.Range("J91").Value = 32460.5
$ws.Range("L91").Value = 32460.5
$ws.Range("N91").Value = -35268.5

$ws.Range("H99").Value = 4004
$ws.Range("I99").Value = 3006
$ws.Range("K99").Value = 3006
$ws.Range("M99").Value = -1508

$ws.Range("H107").Value = 1611.5385
$ws.Range("I107").Value = 1124.75
$ws.Range("K107").Value = 1124.75
$ws.Range("M107").Value = 795.25

$ws.Range("H126").Value = 4004
$ws.Range("I126").Value = 3006
$ws.Range("K126").Value = 9018
$ws.Range("M126").Value = -6548

$ws.Range("H132").Value = 1722
$ws.Range("I132").Value = 1758.6923
$ws.Range("K132").Value = 5276.0769
$ws.Range("M132").Value = -2746.0769

$ws.Range("H134").Value = 504273.3
$ws.Range("J134").Value = 7101.3
$ws.Range("L134").Value = 21303.9
$ws.Range("N134").Value = -26373.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 8000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 8000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 24000
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -28492

$ws.Range("H139").Value = 2529.4119
$ws.Range("J139").Value = 2468.75
$ws.Range("L139").Value = 7406.25
$ws.Range("N139").Value = -17686.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4024.8235
$ws.Range("I113").Value = 3819.5715
$ws.Range("K113").Value = 3819.5715
$ws.Range("M113").Value = -1649.5715

$ws.Range("H132").Value = 58826056
$ws.Range("I132").Value = 62502530
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 187507590
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -187505060
$ws.Range("N132").Value = -12560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H119").Value = 94995
$ws.Range("J119").Value = 94995
$ws.Range("L119").Value = 94995
$ws.Range("N119").Value = -104671

$ws.Range("H122").Value = 6844.3335
$ws.Range("I122").Value = 5962.778
$ws.Range("K122").Value = 17888.334
$ws.Range("M122").Value = -15438.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2997.5
$ws.Range("I81").Value = 3000
$ws.Range("K81").Value = 6000
$ws.Range("M81").Value = -4939

$ws.Range("H84").Value = 2997.5
$ws.Range("I84").Value = 3000
$ws.Range("K84").Value = 30000
$ws.Range("M84").Value = -24696

$ws.Range("H107").Value = 41668256
$ws.Range("J107").Value = 888.5
$ws.Range("L107").Value = 2665.5
$ws.Range("N107").Value = -6505.5

$ws.Range("H136").Value = 10673.077
$ws.Range("I136").Value = 11145.833
$ws.Range("K136").Value = 33437.499
$ws.Range("M136").Value = -30887.499
